$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "(203957296, Omri Ben Shabat: 4,1)"
$ws.Range("B1").Value = "(206532695, Matan Vakrat: -9,0)"
$ws.Range("C1").Value = "(302962915, Asher  Odeh: 0,-5)"
$ws.Range("D1").Value = "(308035542, Anastasia  Kubi: 0,-2)"
$ws.Range("E1").Value = "(311177802, Christina  Uksusman: 9,-1)"
$ws.Range("F1").Value = "(305251175, Or  Leder: -5,0)"
$ws.Range("G1").Value = "(308051846, Eyal  Sofer: -10,-9)"

$ws.Range("A3").Value = "cost: 659.0273893066587"
$ws.Range("A4").Value = "time: 79.25342366333234"
